# Update Leave Card - add monthly EARNED leave credits for 2023 (Feb-Dec),
# a SOLO PARENT leave entry (Oct 2023), a Forced Leave absence (Dec 2023),
# and start a new "2024" year section in the leave card table (Sheet1 / CONVERTION table Table1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# --- Monthly EARNED (1.25) leave credits for Feb 2023 - Nov 2023 (rows 99-108) ---
$ws.Range("C99").Value  = 1.25   # 02/01/2023
$ws.Range("C100").Value = 1.25   # 03/01/2023
$ws.Range("C101").Value = 1.25   # 04/01/2023
$ws.Range("C102").Value = 1.25   # 05/01/2023
$ws.Range("C103").Value = 1.25   # 06/01/2023
$ws.Range("C104").Value = 1.25   # 07/01/2023
$ws.Range("C105").Value = 1.25   # 08/01/2023
$ws.Range("C106").Value = 1.25   # 09/01/2023

# --- 10/01/2023 row: SOLO PARENT leave particulars + remarks ---
$ws.Range("B107").Value = "SOLOP(2-0-0)"
$ws.Range("C107").Value = 1.25
$ws.Range("K107").Value = "10/31 - 11/3/2023"

# --- 11/01/2023 row ---
$ws.Range("C108").Value = 1.25

# --- 12/01/2023 row: Forced Leave particulars + absence + remarks ---
$ws.Range("B109").Value = "FL(3-0-0)"
$ws.Range("D109").Value = 3
$ws.Range("K109").Value = "12/27-29/2023"

# --- Insert a new row for the "2024" year header, pushing existing rows down ---
$ws.Rows.Item(110).Insert()
$lo.Resize($ws.Range("A8:K133"))

# Copy the formatting from the existing "2023" year-header row (row 97) so the
# new row matches the established year-header style (bold, borders, etc.)
$ws.Range("A97:K97").Copy()
$ws.Range("A110:K110").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Year header text (stored as text, not parsed as a date)
$ws.Range("A110").Value = "'2024"

# Restore the calculated column formula for the new row (table calculated column)
$ws.Range("G110").Formula = '=IF(ISBLANK([@EARNED]),"",[@EARNED])'

$wb.Application.Calculate()
